$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 no longer holds "Travalab01" - clear its value but keep its existing style
$ws.Range("B2").Value = $null

# B3 now holds "Travalab01" (previously empty, plain style) - copy the data style from A2
$ws.Range("A2").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B3").Value = "Travalab01"

# A4 now holds "support@travalab.com" and B4 holds "Travalab01" - apply the same data style
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A4").Value = "support@travalab.com"

$ws.Range("A2").Copy()
$ws.Range("B4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B4").Value = "Travalab01"

$excel.CutCopyMode = $false

# Update the active selection to B2
$ws.Range("B2").Select()
